# Login manager Req col
# 1. SREQ005's requirement text wording update: "lobby" -> "game"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "Each player will be assigned a distinct color on entry into the game"

# 2. Append a new requirement row (SREQ042) about the Login Manager enforcing a
#    minimum nickname/handle length, sourced from uta.cse3310.LoginManager.registerUser()
$ws.Range("A42").Value = "SREQ042"
$ws.Range("B42").Value = "F"
$ws.Range("D42").Value = "UI-Login Screen"
$ws.Range("F42").Value = "User unique ""Nick"" or ""Handle"" will be > length 2"

# E42 reuses the same "Implemented by" text/format as the other LoginManager rows
# (copy the formatting from E4, which already holds this text, so the new cell
# picks up the same border/centered style instead of the default unformatted look)
$ws.Range("E42").Value = "uta.cse3310.LoginManager.registerUser()"
$ws.Range("E4").Copy()
$ws.Range("E42").PasteSpecial(-4122)

# 3. Leave the new last row selected, matching where the author ended up editing
$ws.Range("A42").Select()
